$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-68). Bump it from 46060 (2026-02-07) to 46061 (2026-02-08) for
# every row, leaving everything else (formatting, other columns) untouched.
for ($r = 2; $r -le 68; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value2 = 46061
    }
}
